$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44196
$ws.Range("M2").Value = 56

$ws.Range("D3").Value = 44175
$ws.Range("M3").Value = 25

$ws.Range("D4").Value = 44181
$ws.Range("M4").Value = 30

$ws.Range("D6").Value = 44179
$ws.Range("M6").Value = 45
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 20000
$ws.Range("P6").Value = 20000
$ws.Range("S6").Value = 4000

$ws.Range("D7").Value = 44193

$ws.Range("D8").Value = 44189
$ws.Range("M8").Value = 40

$ws.Range("D9").Value = 44188
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 15000
$ws.Range("S9").Value = 3000
